$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 61
$ws.Range("B61").Value = 6905571
$ws.Range("F61").Value = 'FBC Melgar'
$ws.Range("G61").Value = 'Sporting Cristal'
$ws.Range("I61").Value = 1
$ws.Range("J61").Value = 'D'
$ws.Range("K61").Value = 2.1
$ws.Range("L61").Value = 3.4
$ws.Range("M61").Value = 3
$ws.Range("N61").Value = 1.75
$ws.Range("O61").Value = 3.8
$ws.Range("P61").Value = 4.75
$ws.Range("Q61").Value = -0.75
$ws.Range("R61").Value = 1.95
$ws.Range("S61").Value = 1.85
$ws.Range("U61").Value = 1.95
$ws.Range("V61").Value = 1.85
$ws.Range("W61").Value = -1
$ws.Range("X61").Value = 2.8
$ws.Range("Z61").Value = -1
$ws.Range("AA61").Value = 0.8500000000000001
$ws.Range("AC61").Value = 0.8500000000000001

# Row 62
$ws.Range("B62").Value = 6905578
$ws.Range("F62").Value = 'AD Tarma'
$ws.Range("G62").Value = 'Atletico Grau'
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 'H'
$ws.Range("K62").Value = 1.75
$ws.Range("L62").Value = 3.6
$ws.Range("M62").Value = 4
$ws.Range("N62").Value = 1.571
$ws.Range("O62").Value = 4.2
$ws.Range("P62").Value = 5.75
$ws.Range("Q62").Value = -1
$ws.Range("R62").Value = 1.975
$ws.Range("S62").Value = 1.825
$ws.Range("U62").Value = 1.8
$ws.Range("V62").Value = 2
$ws.Range("W62").Value = 0.571
$ws.Range("X62").Value = -1
$ws.Range("Z62").Value = 0
$ws.Range("AA62").Value = -0.0
$ws.Range("AC62").Value = 1

# Row 156
$ws.Range("B156").Value = 7211640
$ws.Range("F156").Value = 'UTC Cajamarca'
$ws.Range("G156").Value = 'Sport Boys'
$ws.Range("H156").Value = 1
$ws.Range("I156").Value = 1
$ws.Range("J156").Value = 'D'
$ws.Range("K156").Value = 1.615
$ws.Range("L156").Value = 3.75
$ws.Range("M156").Value = 5
$ws.Range("N156").Value = 1.5
$ws.Range("O156").Value = 4.2
$ws.Range("P156").Value = 6.5
$ws.Range("Q156").Value = -1
$ws.Range("R156").Value = 1.8
$ws.Range("S156").Value = 2.05
$ws.Range("T156").Value = 2.5
$ws.Range("U156").Value = 1.875
$ws.Range("V156").Value = 1.975
$ws.Range("W156").Value = -1
$ws.Range("X156").Value = 3.2
$ws.Range("Z156").Value = -1
$ws.Range("AA156").Value = 1.05
$ws.Range("AC156").Value = 0.9750000000000001

# Row 157
$ws.Range("B157").Value = 7211641
$ws.Range("F157").Value = 'Sport Huancayo'
$ws.Range("G157").Value = 'Deportivo Municipal'
$ws.Range("H157").Value = 2
$ws.Range("I157").Value = 0
$ws.Range("J157").Value = 'H'
$ws.Range("K157").Value = 1.125
$ws.Range("L157").Value = 7
$ws.Range("M157").Value = 17
$ws.Range("N157").Value = 1.166
$ws.Range("O157").Value = 6.5
$ws.Range("P157").Value = 12
$ws.Range("Q157").Value = -2
$ws.Range("R157").Value = 1.775
$ws.Range("S157").Value = 2.025
$ws.Range("T157").Value = 3.5
$ws.Range("U157").Value = 1.9
$ws.Range("V157").Value = 1.9
$ws.Range("W157").Value = 0.1659999999999999
$ws.Range("X157").Value = -1
$ws.Range("Z157").Value = 0
$ws.Range("AA157").Value = -0.0
$ws.Range("AC157").Value = 0.8999999999999999

# Row 175
$ws.Range("B175").Value = 7302200
$ws.Range("F175").Value = 'Carlos Manucci'
$ws.Range("G175").Value = 'Deportivo Binacional'
$ws.Range("H175").Value = 3
$ws.Range("I175").Value = 2
$ws.Range("K175").Value = 2
$ws.Range("L175").Value = 3.2
$ws.Range("M175").Value = 3.75
$ws.Range("N175").Value = 1.75
$ws.Range("O175").Value = 3.4
$ws.Range("P175").Value = 4.333
$ws.Range("Q175").Value = -0.5
$ws.Range("R175").Value = 1.85
$ws.Range("S175").Value = 1.95
$ws.Range("T175").Value = 2.5
$ws.Range("U175").Value = 1.85
$ws.Range("V175").Value = 1.95
$ws.Range("W175").Value = 0.75
$ws.Range("Z175").Value = 0.8500000000000001
$ws.Range("AA175").Value = -1
$ws.Range("AB175").Value = 0.8500000000000001
$ws.Range("AC175").Value = -1

# Row 176
$ws.Range("B176").Value = 7302795
$ws.Range("F176").Value = 'Unin Comercio'
$ws.Range("G176").Value = 'Deportivo Garcilaso'
$ws.Range("H176").Value = 1
$ws.Range("J176").Value = 'A'
$ws.Range("K176").Value = 2.25
$ws.Range("L176").Value = 3.3
$ws.Range("M176").Value = 2.7
$ws.Range("O176").Value = 3.6
$ws.Range("P176").Value = 4
$ws.Range("R176").Value = 1.8
$ws.Range("S176").Value = 2
$ws.Range("T176").Value = 2.75
$ws.Range("U176").Value = 1.825
$ws.Range("V176").Value = 1.975
$ws.Range("W176").Value = -1
$ws.Range("Y176").Value = 3
$ws.Range("Z176").Value = -1
$ws.Range("AA176").Value = 1
$ws.Range("AB176").Value = 0.4125
$ws.Range("AC176").Value = -0.5

# Row 177
$ws.Range("B177").Value = 7302796
$ws.Range("F177").Value = 'Sport Huancayo'
$ws.Range("G177").Value = 'Sport Boys'
$ws.Range("I177").Value = 0
$ws.Range("J177").Value = 'H'
$ws.Range("K177").Value = 1.727
$ws.Range("L177").Value = 3.75
$ws.Range("M177").Value = 4.333
$ws.Range("N177").Value = 1.25
$ws.Range("O177").Value = 5.25
$ws.Range("P177").Value = 10
$ws.Range("Q177").Value = -1.75
$ws.Range("R177").Value = 1.925
$ws.Range("S177").Value = 1.875
$ws.Range("T177").Value = 3
$ws.Range("U177").Value = 1.875
$ws.Range("V177").Value = 1.925
$ws.Range("W177").Value = 0.25
$ws.Range("Y177").Value = -1
$ws.Range("AA177").Value = 0.875
$ws.Range("AB177").Value = -1
$ws.Range("AC177").Value = 0.925

# Row 184
$ws.Range("B184").Value = 7384626
$ws.Range("F184").Value = 'Sporting Cristal'
$ws.Range("G184").Value = 'Alianza Atletico'
$ws.Range("H184").Value = 3
$ws.Range("I184").Value = 0
$ws.Range("J184").Value = 'H'
$ws.Range("K184").Value = 1.3
$ws.Range("L184").Value = 5
$ws.Range("M184").Value = 9
$ws.Range("N184").Value = 1.166
$ws.Range("O184").Value = 6.5
$ws.Range("P184").Value = 13
$ws.Range("Q184").Value = -2
$ws.Range("R184").Value = 1.85
$ws.Range("S184").Value = 1.95
$ws.Range("T184").Value = 3.25
$ws.Range("U184").Value = 2
$ws.Range("V184").Value = 1.8
$ws.Range("W184").Value = 0.1659999999999999
$ws.Range("Y184").Value = -1
$ws.Range("Z184").Value = 0.8500000000000001
$ws.Range("AA184").Value = -1
$ws.Range("AB184").Value = -0.5
$ws.Range("AC184").Value = 0.4

# Row 185
$ws.Range("B185").Value = 7384628
$ws.Range("F185").Value = 'Deportivo Binacional'
$ws.Range("G185").Value = 'FBC Melgar'
$ws.Range("H185").Value = 1
$ws.Range("I185").Value = 2
$ws.Range("J185").Value = 'A'
$ws.Range("K185").Value = 2.75
$ws.Range("L185").Value = 3.3
$ws.Range("M185").Value = 2.375
$ws.Range("N185").Value = 3.3
$ws.Range("O185").Value = 3.6
$ws.Range("P185").Value = 2
$ws.Range("Q185").Value = 0.5
$ws.Range("R185").Value = 1.8
$ws.Range("S185").Value = 2
$ws.Range("T185").Value = 2.75
$ws.Range("U185").Value = 1.975
$ws.Range("V185").Value = 1.875
$ws.Range("W185").Value = -1
$ws.Range("Y185").Value = 1
$ws.Range("Z185").Value = -1
$ws.Range("AA185").Value = 1
$ws.Range("AB185").Value = 0.4875
$ws.Range("AC185").Value = -0.5

# Row 292
$ws.Range("B292").Value = 8042275
$ws.Range("E292").Value = 45403.54166666666
$ws.Range("F292").Value = 'Sporting Cristal'
$ws.Range("G292").Value = 'Cusco FC'
$ws.Range("K292").Value = 1.4
$ws.Range("L292").Value = 4.5
$ws.Range("M292").Value = 8
$ws.Range("N292").Value = 1.25
$ws.Range("O292").Value = 5.75
$ws.Range("P292").Value = 13
$ws.Range("Q292").Value = -1.75
$ws.Range("R292").Value = 1.95
$ws.Range("S292").Value = 1.9
$ws.Range("T292").Value = 3
$ws.Range("U292").Value = 1.825
$ws.Range("V292").Value = 2.025

# Row 293
$ws.Range("B293").Value = 8042082
$ws.Range("E293").Value = 45403.6875
$ws.Range("F293").Value = 'Atletico Grau'
$ws.Range("G293").Value = 'FBC Melgar'
$ws.Range("K293").Value = 2.625
$ws.Range("L293").Value = 3.2
$ws.Range("M293").Value = 2.625
$ws.Range("N293").Value = 2.375
$ws.Range("O293").Value = 3
$ws.Range("P293").Value = 3.1
$ws.Range("Q293").Value = -0.25
$ws.Range("R293").Value = 2.05
$ws.Range("S293").Value = 1.8
$ws.Range("T293").Value = 2
$ws.Range("U293").Value = 1.775
$ws.Range("V293").Value = 2.1

# Row 294
$ws.Range("B294").Value = 8042276
$ws.Range("E294").Value = 45404.70833333334
$ws.Range("F294").Value = 'Comerciantes Unidos'
$ws.Range("G294").Value = 'Alianza Atletico'
$ws.Range("K294").Value = 2.1
$ws.Range("L294").Value = 3.4
$ws.Range("M294").Value = 3.4
$ws.Range("N294").Value = 1.666
$ws.Range("O294").Value = 3.8
$ws.Range("P294").Value = 4.75
$ws.Range("Q294").Value = -0.75
$ws.Range("R294").Value = 1.925
$ws.Range("S294").Value = 1.925
$ws.Range("T294").Value = 2.5
$ws.Range("U294").Value = 1.875
$ws.Range("V294").Value = 1.975
